# Slide 4: "Rounded Rectangle 5" (inside the "Group 13" group) currently reads
#   Run1 (not bold): "Time left for current "
#   Run2 (bold):      "Card"
# Target:
#   Run1 (not bold): "T"
#   Run2 (bold):      "ime"
#   Run3 (not bold): " left for current Card"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$g = $s.Shapes.Item(1)
$shp = $g.GroupItems.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Step 1: grow the first (non-bold) run so it also carries a trailing
# copy of "Card" - this keeps everything inside the ORIGINAL non-bold run,
# so its run properties (no explicit b="" attribute) stay untouched.
$firstRun = $tr.Characters(1, 22)              # "Time left for current "
$firstRun.Text = "T left for current Card"

# --- Step 2: re-expand that same (still non-bold) run so the full phrase
# "Time left for current Card" precedes the original bold "Card" run.
$firstRunGrown = $tr.Characters(1, 23)         # "T left for current Card"
$firstRunGrown.Text = "Time left for current Card"

# --- Step 3: carve "ime" back out of the non-bold run and mark it bold;
# this splits the run into "T" / "ime"(bold) / " left for current Card",
# without ever setting Bold=False on anything, so no stray b="0" appears.
$boldPart = $tr.Characters(2, 3)               # "ime"
$boldPart.Text = "ime"
$boldPart = $tr.Characters(2, 3)
$boldPart.Font.Bold = $true

# --- Step 4: drop the now-duplicated trailing "Card" that still belongs to
# the original bold run.
$dupTail = $tr.Characters(27, 4)               # trailing duplicate "Card"
$dupTail.Text = ""
